$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix 11 cells that flip between numeric value and the "NaN" shared string (index 68) ---
$ws.Range("DB29").Value = 1
$ws.Range("DK53").Value = "NaN"
$ws.Range("DB85").Value = "NaN"
$ws.Range("AD106").Value = "NaN"
$ws.Range("J114").Value = "NaN"
$ws.Range("DW117").Value = "NaN"
$ws.Range("AD119").Value = "NaN"
$ws.Range("AI132").Value = "NaN"
$ws.Range("AI133").Value = "NaN"
$ws.Range("AI140").Value = "NaN"
$ws.Range("AI141").Value = "NaN"

# --- 2) Append a new data row 195 (2020-09-15) with the full A:DX record ---
$row195 = New-Object 'object[,]' 1,128
$row195[0,0] = 44089
$row195[0,1] = 728590
$row195[0,2] = 2725
$row195[0,3] = 97642
$row195[0,4] = 66038
$row195[0,5] = 243456
$row195[0,6] = 27551
$row195[0,7] = 5319
$row195[0,8] = 4191
$row195[0,9] = 7346
$row195[0,10] = 7439
$row195[0,11] = 16113
$row195[0,12] = 3869
$row195[0,13] = 22236
$row195[0,14] = 28918
$row195[0,15] = 6864
$row195[0,16] = 7945
$row195[0,17] = 13970
$row195[0,18] = 12054
$row195[0,19] = 16316
$row195[0,20] = 13837
$row195[0,21] = 3401
$row195[0,22] = 2115
$row195[0,23] = 8428
$row195[0,24] = 25392
$row195[0,25] = 13293
$row195[0,26] = 9865
$row195[0,27] = 54395
$row195[0,28] = 1597
$row195[0,29] = 596
$row195[0,30] = 592
$row195[0,31] = 462
$row195[0,32] = 459
$row195[0,33] = 314
$row195[0,34] = 487
$row195[0,35] = 2007
$row195[0,36] = 4520
$row195[0,37] = 37235
$row195[0,38] = 8530
$row195[0,39] = 2514
$row195[0,40] = 42427
$row195[0,41] = 1052
$row195[0,42] = 22157
$row195[0,43] = 1499
$row195[0,44] = 9495
$row195[0,45] = 1621
$row195[0,46] = 1593
$row195[0,47] = 6603
$row195[0,48] = 1797
$row195[0,49] = 954
$row195[0,50] = 2484
$row195[0,51] = 2654
$row195[0,52] = 56642
$row195[0,53] = 13388
$row195[0,54] = 4670
$row195[0,55] = 8798
$row195[0,56] = 5616
$row195[0,57] = 281
$row195[0,58] = 1438
$row195[0,59] = 2657
$row195[0,60] = 737
$row195[0,61] = 2124
$row195[0,62] = 9268
$row195[0,63] = 9227
$row195[0,64] = 9887
$row195[0,65] = 14140
$row195[0,66] = 1936
$row195[0,67] = 887
$row195[0,68] = 11902
$row195[0,69] = 9562
$row195[0,70] = 11184
$row195[0,71] = 2204
$row195[0,72] = 1883
$row195[0,73] = 4610
$row195[0,74] = 4271
$row195[0,75] = 1513
$row195[0,76] = 5445
$row195[0,77] = 3138
$row195[0,78] = 1784
$row195[0,79] = 865
$row195[0,80] = 2646
$row195[0,81] = 2151
$row195[0,82] = 1703
$row195[0,83] = 1330
$row195[0,84] = 5953
$row195[0,85] = 1884
$row195[0,86] = 1342
$row195[0,87] = 1606
$row195[0,88] = 1961
$row195[0,89] = 1941
$row195[0,90] = 2305
$row195[0,91] = 1439
$row195[0,92] = 1188
$row195[0,93] = 1168
$row195[0,94] = 800
$row195[0,95] = 3209
$row195[0,96] = 1314
$row195[0,97] = 890
$row195[0,98] = 932
$row195[0,99] = 1638
$row195[0,100] = 1463
$row195[0,101] = 721
$row195[0,102] = 830
$row195[0,103] = 1186
$row195[0,104] = 1450
$row195[0,105] = 1322
$row195[0,106] = 1370
$row195[0,107] = 1088
$row195[0,108] = 333
$row195[0,109] = 352
$row195[0,110] = 774
$row195[0,111] = 705
$row195[0,112] = 459
$row195[0,113] = 536
$row195[0,114] = 365
$row195[0,115] = 653
$row195[0,116] = 741
$row195[0,117] = 521
$row195[0,118] = 485
$row195[0,119] = 372
$row195[0,120] = 518
$row195[0,121] = 130536
$row195[0,122] = 309278
$row195[0,123] = 15007
$row195[0,124] = 133161
$row195[0,125] = 82019
$row195[0,126] = 39278
$row195[0,127] = 11206
$ws.Range("A195:DX195").Value = $row195

# --- 3) Match the saved selection state (active cell on the newly added totals) ---
$ws.Range("DR195:DX195").Select() | Out-Null

